$wb = $excel.ActiveWorkbook

# Insert a new "Name" column (column B) on every sheet, shifting the
# existing columns one place to the right, and fill in the new header /
# description cells.

$sheetNames = @("Drilling", "Water Strike", "Stratigraphic Log", "Construction")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns("B").Insert()
    $ws.Range("B1").Value = "Name "
    $ws.Range("B2").Value = "The name of the data point."
    $ws.Columns("B").ColumnWidth = 18.75
}

# "Drilling" is now the active / selected tab (was "Stratigraphic Log").
$drilling = $wb.Worksheets.Item("Drilling")
$drilling.Activate()

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Select()
}

$drilling.Range("B1").Select()
